$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3:D3").Value = 4
$ws.Range("C5:D5").Value = 934
$ws.Range("C7:D7").Value = 292
$ws.Range("C9:D9").Value = 191
$ws.Range("C11:D11").Value = 1051
$ws.Range("C13:D13").Value = 655
$ws.Range("C16:D16").Value = 14
$ws.Range("C17:D17").Value = 457
$ws.Range("C19:D19").Value = 291
$ws.Range("C21:D21").Value = 176
$ws.Range("C23:D23").Value = 828
$ws.Range("C25:D25").Value = 9
$ws.Range("C27:D27").Value = 343
$ws.Range("C29:D29").Value = 612
$ws.Range("C32:D32").Value = 84
$ws.Range("C33:D33").Value = 189
$ws.Range("C35:D35").Value = 1400
$ws.Range("C37:D37").Value = 138
$ws.Range("C39:D39").Value = 2764
$ws.Range("C41:D41").Value = 713
$ws.Range("C42:D42").Value = 519
$ws.Range("C44:D44").Value = 73
$ws.Range("C46:D46").Value = 1985
$ws.Range("C48:D48").Value = 174
$ws.Range("C50:D50").Value = 1623
$ws.Range("C52:D52").Value = 135
$ws.Range("C54:D54").Value = 185
$ws.Range("C56:D56").Value = 937
$ws.Range("C58:D58").Value = 4236
$ws.Range("C60:D60").Value = 67
$ws.Range("C62:D62").Value = 4186
$ws.Range("C64:D64").Value = 555
$ws.Range("C66:D66").Value = 206
$ws.Range("C68:D68").Value = 85
$ws.Range("C70:D70").Value = 1746
$ws.Range("C72:D72").Value = 137
$ws.Range("C74:D74").Value = 234
$ws.Range("C76:D76").Value = 104
$ws.Range("C78:D78").Value = 1386
$ws.Range("C80:D80").Value = 187
$ws.Range("C82:D82").Value = 101
$ws.Range("C84:D84").Value = 95
$ws.Range("C86:D86").Value = 86
$ws.Range("C88:D88").Value = 331
$ws.Range("C90:D90").Value = 769
$ws.Range("C92:D92").Value = 1265
$ws.Range("C94:D94").Value = 714
$ws.Range("C96:D96").Value = 96
$ws.Range("C98:D98").Value = 499
$ws.Range("C100:D100").Value = 204
$ws.Range("C102:D102").Value = 350
$ws.Range("C104:D104").Value = 34
$ws.Range("C106:D106").Value = 345
$ws.Range("C108:D108").Value = 90
$ws.Range("C110:D110").Value = 109
$ws.Range("C112:D112").Value = 186
$ws.Range("C114:D114").Value = 395
$ws.Range("C116:D116").Value = 4373
$ws.Range("C118:D118").Value = 778
$ws.Range("C120:D120").Value = 684
$ws.Range("C122:D122").Value = 158
$ws.Range("C124:D124").Value = 380
$ws.Range("C126:D126").Value = 93
$ws.Range("C128:D128").Value = 94
$ws.Range("C130:D130").Value = 37
$ws.Range("C132:D132").Value = 91
$ws.Range("C134:D134").Value = 340
$ws.Range("C136:D136").Value = 2541
$ws.Range("C138:D138").Value = 159
$ws.Range("C140:D140").Value = 595
$ws.Range("C142:D142").Value = 1703
$ws.Range("C144:D144").Value = 1704
$ws.Range("C145").Value = 684.7777777777778
